# systemConfig.xlsx update: rename the MSE-* comparison sheets and refresh
# their group-arrival / distributed-precoding legends & params.
#
#   MSE-NM20-0   -> MSE-C   (headers become OLC / OLCM variants)
#   MSE-NM20-20  -> MSE-B   (headers become OLM variant on the last 3 cols)
#   MSE-WM20-0   -> MSE-L   (headers become OLE / OLEM variants) -- rebuilt
#                            as a fresh sheet copy so it gets a new sheetId
#   MSE-WM20-20  -> removed entirely

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Build the replacement "MSE-L" sheet from a copy of "MSE-WM20-0" (so the
#    new sheet picks up a brand-new sheetId instead of inheriting the old
#    one), then drop the two sheets that are going away.
# ---------------------------------------------------------------------------
$srcL = $wb.Worksheets("MSE-WM20-0")
$srcL.Copy($wb.Worksheets($wb.Worksheets.Count))
$wb.Worksheets("MSE-WM20-0 (2)").Name = "MSE-L"

$wb.Worksheets("MSE-WM20-0").Delete()
$wb.Worksheets("MSE-WM20-20").Delete()

# ---------------------------------------------------------------------------
# 2. Rename the two sheets that survive in place.
# ---------------------------------------------------------------------------
$wb.Worksheets("MSE-NM20-0").Name = "MSE-C"
$wb.Worksheets("MSE-NM20-20").Name = "MSE-B"

# ---------------------------------------------------------------------------
# 3. MSE-C content updates.
# ---------------------------------------------------------------------------
$wsC = $wb.Worksheets("MSE-C")
$wsC.Range("B1").Value = "OLC (OTA-3,INT-[20,20])"
$wsC.Range("C1").Value = "OLC (OTA-5,INT-[20,20])"
$wsC.Range("D1").Value = "OLC (OTA-10,INT-[20,20])"
$wsC.Range("E1").Value = "OLCM (OTA-3,INT-[20,20])"
$wsC.Range("F1").Value = "OLCM (OTA-5,INT-[20,20])"
$wsC.Range("G1").Value = "OLCM (OTA-10,INT-[20,20])"

$wsC.Range("B6").Value = "Optimal"
$wsC.Range("C6").Value = "Optimal"
$wsC.Range("D6").Value = "Optimal"

$wsC.Range("E8").Value = 10
$wsC.Range("F8").Value = 10
$wsC.Range("G8").Value = 10

$wsC.Range("B9").Value = 20
$wsC.Range("C9").Value = 20
$wsC.Range("D9").Value = 20
$wsC.Range("E9").Value = 20
$wsC.Range("F9").Value = 20
$wsC.Range("G9").Value = 20

# ---------------------------------------------------------------------------
# 4. MSE-B content updates.
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets("MSE-B")
$wsB.Range("E1").Value = "OLM (OTA-3,INT-[20,20])"
$wsB.Range("F1").Value = "OLM (OTA-5,INT-[20,20])"
$wsB.Range("G1").Value = "OLM (OTA-10,INT-[20,20])"

$wsB.Range("E6").Value = "MMSE"
$wsB.Range("F6").Value = "MMSE"
$wsB.Range("G6").Value = "MMSE"

$wsB.Range("E8").Value = 10
$wsB.Range("F8").Value = 10
$wsB.Range("G8").Value = 10

$wsB.Range("B9").Value = 20
$wsB.Range("C9").Value = 20
$wsB.Range("D9").Value = 20
$wsB.Range("E9").Value = 20
$wsB.Range("F9").Value = 20
$wsB.Range("G9").Value = 20

# ---------------------------------------------------------------------------
# 5. MSE-L content updates.
# ---------------------------------------------------------------------------
$wsL = $wb.Worksheets("MSE-L")
$wsL.Range("B1").Value = "OLE (OTA-3,INT-[20,20])"
$wsL.Range("C1").Value = "OLE (OTA-5,INT-[20,20])"
$wsL.Range("D1").Value = "OLE (OTA-10,INT-[20,20])"
$wsL.Range("E1").Value = "OLEM (OTA-3,INT-[20,20])"
$wsL.Range("F1").Value = "OLEM (OTA-5,INT-[20,20])"
$wsL.Range("G1").Value = "OLEM (OTA-10,INT-[20,20])"

$wsL.Range("B6").Value = "E-Optimal"
$wsL.Range("C6").Value = "E-Optimal"
$wsL.Range("D6").Value = "E-Optimal"
$wsL.Range("E6").Value = "E-Optimal"
$wsL.Range("F6").Value = "E-Optimal"
$wsL.Range("G6").Value = "E-Optimal"

$wsL.Range("B8").Value = 1
$wsL.Range("C8").Value = 1
$wsL.Range("D8").Value = 1

$wsL.Range("B9").Value = 20
$wsL.Range("C9").Value = 20
$wsL.Range("D9").Value = 20
$wsL.Range("E9").Value = 20
$wsL.Range("F9").Value = 20
$wsL.Range("G9").Value = 20

# ---------------------------------------------------------------------------
# 6. View-state: per-sheet selection, then leave MSE-L as the active tab
#    (selecting last is what makes a sheet the active one).
# ---------------------------------------------------------------------------
$wb.Worksheets("MSE-O").Range("D26").Select()
$wsC.Range("D33").Select()
$wsB.Range("C31").Select()
$wsL.Range("C34").Select()
